$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Update Status column values for specific backlog items
$ws.Range("F6").Value = "Done"
$ws.Range("F10").Value = "Doing"
$ws.Range("F11").Value = "Done"

# Update the selected cell on the sheet
$ws.Range("D13").Select()
